# Add FreeRTOS v9.0.0 source code entry to the project directory list.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 8 with the new directory entry (mirrors existing rows 4-7).
$ws.Range("B8").Value = "FreeRTOSv9.0.0"
$ws.Range("C8").Value = "FreeRTOS 9.0.0 源码"

# Move the active selection, matching the recorded cursor position after edit.
$ws.Range("E23").Select()
